$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.557.23"
$ws.Range("E2").Value = "  +1.63%  "
$ws.Range("D3").Value = "2.297.51"
$ws.Range("E3").Value = "  +0.95%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'156.12"
$ws.Range("E5").Value = "  +15,495.04%  "
$ws.Range("D6").Value = "'307.56"
$ws.Range("E6").Value = "  +0.66%  "
$ws.Range("D7").Value = "'96.84"
$ws.Range("E7").Value = "  +4.51%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("E10").Value = "  +2.55%  "
$ws.Range("D11").Value = "'35.72"
$ws.Range("E11").Value = "  +8.82%  "
$ws.Range("D12").Value = "'0.0813"
$ws.Range("E12").Value = "  +1.62%  "
$ws.Range("E13").Value = "  -1.61%  "
$ws.Range("D14").Value = "'6.76"
$ws.Range("E14").Value = "  +1.73%  "
$ws.Range("D15").Value = "2.651.60"
$ws.Range("E15").Value = "  +0.89%  "
$ws.Range("E16").Value = "  +1.97%  "
$ws.Range("D17").Value = "2.309.34"
$ws.Range("E17").Value = "  +0.43%  "
$ws.Range("D18").Value = "'0.799"
$ws.Range("E18").Value = "  +4.42%  "
$ws.Range("D19").Value = "42.405.18"
$ws.Range("D20").Value = "'12.88"
$ws.Range("E20").Value = "  +4.18%  "
$ws.Range("E21").Value = "  +1.52%  "
$ws.Range("D22").Value = "'6.06"
$ws.Range("E22").Value = "  +1.89%  "
$ws.Range("D23").Value = "'68.29"
$ws.Range("E23").Value = "  +1.72%  "
$ws.Range("D24").Value = "'245.45"
$ws.Range("E24").Value = "  +1.05%  "
$ws.Range("E25").Value = "  +0.97%  "
$ws.Range("E26").Value = "  +1.65%  "
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("D28").Value = "'24.32"
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").Value = "'36.83"
$ws.Range("E29").Value = "  +7.47%  "
$ws.Range("D30").Value = "'9.74"
$ws.Range("E30").Value = "  +1.25%  "
$ws.Range("E31").Value = "  +1.51%  "
$ws.Range("D32").Value = "'161.46"
$ws.Range("E32").Value = "  +1.56%  "
$ws.Range("D33").Value = "'5.40"
$ws.Range("E33").Value = "  +4.01%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("D35").Value = "'0.0757"
$ws.Range("E35").Value = "  +0.79%  "
$ws.Range("D36").Value = "'3.12"
$ws.Range("E36").Value = "  +2.37%  "
$ws.Range("D37").Value = "'17.54"
$ws.Range("E37").Value = "  +2.56%  "
$ws.Range("E38").Value = "  +4.93%  "
$ws.Range("E39").Value = "  +0.33%  "
$ws.Range("D40").Value = "'1.86"
$ws.Range("E40").Value = "  +1.76%  "
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("D42").Value = "'4.22"
$ws.Range("E42").Value = "  +7.35%  "
$ws.Range("D43").Value = "'20.17"
$ws.Range("E43").Value = "  +2.89%  "
$ws.Range("D44").Value = "2.020.16"
$ws.Range("E44").Value = "  -2.46%  "
$ws.Range("E45").Value = "  +10.68%  "
$ws.Range("E46").Value = "  +2.51%  "
$ws.Range("D47").Value = "'10.31"
$ws.Range("E47").Value = "  -0.48%  "
$ws.Range("D48").Value = "'3.02"
$ws.Range("E48").Value = "  +3.89%  "
$ws.Range("D49").Value = "'53.99"
$ws.Range("E49").Value = "  +4.16%  "
$ws.Range("D50").Value = "'1.56"
$ws.Range("E50").Value = "  +2.14%  "
$ws.Range("D51").Value = "'73.21"
$ws.Range("E51").Value = "  +0.21%  "
